$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Explicit black font color for the corner header cell (A1), replacing the
# theme-based black it inherited before.
$ws.Range("A1").Font.Color = 0

# Tighten the data-row heights slightly (uniform -0.75pt shift), leaving the
# taller wrapped row (5) untouched.
$ws.Rows.Item(1).RowHeight = 20.25
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
$ws.Rows.Item(4).RowHeight = 19.5
$ws.Rows.Item(6).RowHeight = 19.5
$ws.Rows.Item(7).RowHeight = 19.5
$ws.Rows.Item(8).RowHeight = 19.5
$ws.Rows.Item(9).RowHeight = 19.5
$ws.Rows.Item(10).RowHeight = 19.5
$ws.Rows.Item(11).RowHeight = 20.25
